# Weekly update for "Fruta, Feria Lagunitas de Puerto Montt - Kiwi":
# three new daily price records (kiwi Hayward, Especial/Primera/Segunda
# quality grades) are inserted at the top of the existing data block
# (current rows 513-515), pushing the rest of the table down by three
# rows (old row 513 -> new row 516, ..., old row 533 -> new row 536).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right above the current row 513.
$ws.Rows("513:515").Insert()

# Common columns shared by every record in this block.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad  = "Hayward"
$unidad    = "`$/caja 15 kilos"
$origen    = "Región de O'Higgins"
$kgUnidad  = 15

$newRows = @(
    @{ Row = 513; Fecha = 45075; Calidad = "Especial"; Volumen = 100; PMin = 19000; PMax = 19000; PProm = 19000; PrecioKg = 1267 },
    @{ Row = 514; Fecha = 45075; Calidad = "Primera";  Volumen = 100; PMin = 17000; PMax = 17000; PProm = 17000; PrecioKg = 1133 },
    @{ Row = 515; Fecha = 45075; Calidad = "Segunda";  Volumen = 100; PMin = 15000; PMax = 15000; PProm = 15000; PrecioKg = 1000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
